$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3860
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 4075
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 4075
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -4213

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2218.7273
$ws.Range("I135").Value = 1458.6666
$ws.Range("K135").Value = 13127.9994
$ws.Range("M135").Value = -10592.9994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644.4878
$ws.Range("I2").Value = 539.69446
$ws.Range("K2").Value = 539.69446
$ws.Range("M2").Value = -426.69446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28386.916
$ws.Range("I32").Value = 28870.857
$ws.Range("J32").Value = 24999.334
$ws.Range("K32").Value = 28870.857
$ws.Range("L32").Value = 24999.334
$ws.Range("M32").Value = -28583.857
$ws.Range("N32").Value = -25573.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 64013.875
$ws.Range("I74").Value = 78262.16
$ws.Range("K74").Value = 78262.16
$ws.Range("M74").Value = -77388.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 64013.875
$ws.Range("I77").Value = 78262.16
$ws.Range("K77").Value = 391310.8
$ws.Range("M77").Value = -386942.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 14226.333
$ws.Range("I102").Value = 22642.6
$ws.Range("J102").Value = 3706
$ws.Range("K102").Value = 22642.6
$ws.Range("L102").Value = 3706
$ws.Range("M102").Value = -21020.6
$ws.Range("N102").Value = -6950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 644.4878
$ws.Range("I116").Value = 539.69446
$ws.Range("K116").Value = 539.69446
$ws.Range("M116").Value = 1754.30554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 24579.422
$ws.Range("I132").Value = 26219.244
$ws.Range("K132").Value = 78657.73199999999
$ws.Range("M132").Value = -76127.73199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644.4878
$ws.Range("I3").Value = 539.69446
$ws.Range("K3").Value = 539.69446
$ws.Range("M3").Value = -425.69446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 23001
$ws.Range("J6").Value = 23001
$ws.Range("L6").Value = 23001
$ws.Range("N6").Value = -23227

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2228.3635
$ws.Range("I20").Value = 1701.5
$ws.Range("K20").Value = 1701.5
$ws.Range("M20").Value = -1454.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4487
$ws.Range("I105").Value = 4330.4614
$ws.Range("K105").Value = 4330.4614
$ws.Range("M105").Value = -2583.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2434.762
$ws.Range("I134").Value = 1948.9412
$ws.Range("K134").Value = 5846.8236
$ws.Range("M134").Value = -3311.8236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4061.7307
$ws.Range("I31").Value = 2146.0557
$ws.Range("J31").Value = 8372
$ws.Range("K31").Value = 2146.0557
$ws.Range("L31").Value = 8372
$ws.Range("M31").Value = -1851.0557
$ws.Range("N31").Value = -8962

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4061.7307
$ws.Range("I34").Value = 2146.0557
$ws.Range("J34").Value = 8372
$ws.Range("K34").Value = 2146.0557
$ws.Range("L34").Value = 8372
$ws.Range("M34").Value = -1944.0557
$ws.Range("N34").Value = -8776

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 53138.2
$ws.Range("I58").Value = 58740.89
$ws.Range("J58").Value = 2714
$ws.Range("K58").Value = 58740.89
$ws.Range("L58").Value = 2714
$ws.Range("M58").Value = -58537.89
$ws.Range("N58").Value = -3120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 55000
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 447.18182
$ws.Range("J107").Value = 150
$ws.Range("L107").Value = 150
$ws.Range("N107").Value = -3990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 22000
$ws.Range("J131").Value = 22000
$ws.Range("L131").Value = 22000
$ws.Range("N131").Value = -32080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5235
$ws.Range("I132").Value = 5235
$ws.Range("K132").Value = 15705
$ws.Range("M132").Value = -13175

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 53138.2
$ws.Range("I136").Value = 58740.89
$ws.Range("J136").Value = 2714
$ws.Range("K136").Value = 176222.67
$ws.Range("L136").Value = 8142
$ws.Range("M136").Value = -173672.67
$ws.Range("N136").Value = -13242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2299
$ws.Range("I3").Value = 2299
$ws.Range("K3").Value = 6897
$ws.Range("M3").Value = -6785

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 782.7
$ws.Range("J5").Value = 767.2
$ws.Range("L5").Value = 2301.6
$ws.Range("N5").Value = -2525.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 72940.42999999999
$ws.Range("J37").Value = 72940.42999999999
$ws.Range("L37").Value = 218821.29
$ws.Range("N37").Value = -219045.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3699
$ws.Range("J63").Value = 3699
$ws.Range("L63").Value = 11097
$ws.Range("N63").Value = -12595

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3699
$ws.Range("J66").Value = 3699
$ws.Range("L66").Value = 33291
$ws.Range("N66").Value = -40779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1149.7097
$ws.Range("I113").Value = 794.3333
$ws.Range("J113").Value = 1295.091
$ws.Range("K113").Value = 2382.9999
$ws.Range("L113").Value = 3885.273
$ws.Range("M113").Value = -212.9998999999998
$ws.Range("N113").Value = -8225.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 8914.454
$ws.Range("I126").Value = 3562.2222
$ws.Range("K126").Value = 10686.6666
$ws.Range("M126").Value = -5746.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 894.9167
$ws.Range("I132").Value = 656.4286
$ws.Range("J132").Value = 1228.8
$ws.Range("K132").Value = 5907.8574
$ws.Range("L132").Value = 11059.2
$ws.Range("M132").Value = -3377.8574
$ws.Range("N132").Value = -16119.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 6580.4546
$ws.Range("J133").Value = 8642.857
$ws.Range("L133").Value = 25928.571
$ws.Range("N133").Value = -36048.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 782.7
$ws.Range("J135").Value = 767.2
$ws.Range("L135").Value = 6904.8
$ws.Range("N135").Value = -11974.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 30302
$ws.Range("I46").Value = 19313.666
$ws.Range("J46").Value = 46784.5
$ws.Range("K46").Value = 19313.666
$ws.Range("L46").Value = 46784.5
$ws.Range("M46").Value = -19157.666
$ws.Range("N46").Value = -47096.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 76874.25
$ws.Range("J105").Value = 76874.25
$ws.Range("L105").Value = 76874.25
$ws.Range("N105").Value = -83862.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3334.6
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5450.4287
$ws.Range("I126").Value = 4765.7827
$ws.Range("K126").Value = 14297.3481
$ws.Range("M126").Value = -11827.3481

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7307
$ws.Range("I7").Value = 6568.7144
$ws.Range("J7").Value = 8599
$ws.Range("K7").Value = 6568.7144
$ws.Range("L7").Value = 8599
$ws.Range("M7").Value = -6456.7144
$ws.Range("N7").Value = -8823

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1585.3077
$ws.Range("I40").Value = 1585.3077
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1585.3077
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1449.3077
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3725.1304
$ws.Range("I122").Value = 2971.5217
$ws.Range("J122").Value = 4478.7393
$ws.Range("K122").Value = 8914.5651
$ws.Range("L122").Value = 13436.2179
$ws.Range("M122").Value = -6464.5651
$ws.Range("N122").Value = -18336.2179

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7307
$ws.Range("I126").Value = 6568.7144
$ws.Range("J126").Value = 8599
$ws.Range("K126").Value = 19706.1432
$ws.Range("L126").Value = 25797
$ws.Range("M126").Value = -17236.1432
$ws.Range("N126").Value = -30737

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 67996.25
$ws.Range("J133").Value = 67996.25
$ws.Range("L133").Value = 67996.25
$ws.Range("N133").Value = -73056.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5337.9165
$ws.Range("J136").Value = 6150.5713
$ws.Range("L136").Value = 18451.7139
$ws.Range("N136").Value = -23551.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3941.6667
$ws.Range("J81").Value = 8731.666999999999
$ws.Range("L81").Value = 17463.334
$ws.Range("N81").Value = -19585.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3941.6667
$ws.Range("J84").Value = 8731.666999999999
$ws.Range("L84").Value = 87316.67
$ws.Range("N84").Value = -97924.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 99188
$ws.Range("J109").Value = 99188
$ws.Range("L109").Value = 99188
$ws.Range("N109").Value = -101962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 305615
$ws.Range("I132").Value = 385799.5
$ws.Range("K132").Value = 1157398.5
$ws.Range("M132").Value = -1154868.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3843.6562
$ws.Range("I136").Value = 4227.364
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 12682.092
$ws.Range("M136").Value = -10132.092
